$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.734.45'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.602.69'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.38'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '28.37'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +5.52%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.255'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.85%  '
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '1.831.83'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = '1.603.17'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('D15').Value = '29.713.39'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.78'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.99'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '242.12'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.98'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.47%  '
$ws.Range('D20').Value = '0.0₃0698'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.43'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.11'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '155.29'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0478'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D34').Value = '1.420.19'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.87'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.545'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '55.62'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0494'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.87%  '
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.998'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +19.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '67.00'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.49%  '
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('D49').Value = '1.741.65'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '86.56'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0105'
$ws.Range('E51').Value = '  +4.70%  '
